$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.506.06"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.824.08"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'316.84"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.5170"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").Value = "'0.3892"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "'0.08521"
$ws.Range("E9").Value = "  +9.94%  "
$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'6.443"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "'7.535"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "1.815.77"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'0.00001141"
$ws.Range("E17").Value = "  +4.90%  "
$ws.Range("D18").Value = "'92.93"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'0.06611"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'17.77"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "28.528.21"
$ws.Range("D24").Value = "'11.46"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").Value = "'2.278"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").Value = "'21.03"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "2.025.19"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").Value = "'2.403"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "'125.63"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'0.1090"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("D32").Value = "'1.100"
$ws.Range("E32").Value = "  -4.21%  "
$ws.Range("D33").Value = "'5.731"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'0.07429"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "'3.652"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'0.2234"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'0.02354"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'5.226"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "'8.839"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'0.6325"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "'11.41"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'1.196"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'13.49"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'3.787"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").Value = "'0.5958"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "'126.31"
$ws.Range("D48").Value = "'1.993"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").Value = "'0.06977"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "'74.42"
$ws.Range("E51").Value = "  -0.13%  "
